# Actualización de horarios Línea 141 - 367
# Se agregan nuevas filas scrapeadas a las hojas LP1912, LP1912-215 y 6203-6173
# con el timestamp de actualización 05:19:24.

$wb = $excel.ActiveWorkbook

$nuevaHora = "05:19:24"

# ---------------------------------------------------------------
# Hoja 1: LP1912  (24 -> 29 filas de datos; dimension A1:E29 -> A1:E34)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $nuevaHora"
$ws1.Range("A3").Value = "Total filas: 29"

$ws1Rows = @(
    @($nuevaHora, "06:59", "14_ABASTO",      100, "LP1912"),
    @($nuevaHora, "07:05", "15_ABASTO",      106, "LP1912"),
    @($nuevaHora, "07:07", "225_GOMEZ",      108, "LP1912"),
    @($nuevaHora, "07:11", "215A_EL PATO",   112, "LP1912"),
    @($nuevaHora, "07:15", "11_ETCHEVERRY",  116, "LP1912")
)

$r = 30
foreach ($row in $ws1Rows) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------
# Hoja 2: LP1912-215  (6 -> 7 filas de datos; dimension A1:E11 -> A1:E12)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $nuevaHora"
$ws2.Range("A3").Value = "Total filas: 7"

$ws2.Cells.Item(12, 1).Value = $nuevaHora
$ws2.Cells.Item(12, 2).Value = "07:11"
$ws2.Cells.Item(12, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(12, 4).Value = 112
$ws2.Cells.Item(12, 5).Value = "LP1912"

# ---------------------------------------------------------------
# Hoja 3: 6203-6173  (6 -> 7 filas de datos; dimension A1:E11 -> A1:E12)
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $nuevaHora"
$ws3.Range("A3").Value = "Total filas: 7"

$ws3.Cells.Item(12, 1).Value = $nuevaHora
$ws3.Cells.Item(12, 2).Value = "07:00"
$ws3.Cells.Item(12, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws3.Cells.Item(12, 4).Value = 101
$ws3.Cells.Item(12, 5).Value = "L6173"
